$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# New RMA batch "MO6T" replacing the previous "YTTY" values in row 2-4
$ws.Range("E2").Value = "RMA-MO6T-001"
$ws.Range("F2").Value = "RMA-MO6T-1-1"
$ws.Range("J2").Value = "a7s5f000000xK6kAAE"

$ws.Range("E3").Value = "RMA-MO6T-002"
$ws.Range("F3").Value = "RMA-MO6T-1-2"
$ws.Range("J3").Value = "a7s5f000000xK6lAAE"

$ws.Range("E4").Value = "RMA-MO6T-003"
$ws.Range("F4").Value = "RMA-MO6T-1-3"
$ws.Range("J4").Value = "a7s5f000000xK6mAAE"

Write-Host "done"
